$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 1")
$tr = $sh.TextFrame.TextRange

$oldText = "Clustering using the centroid of hotel locations gave us 130 centers to create Voronoi cells with a reasonable distribution of listings. 36% of Airbnb listings were disregarded for being too far away from a center."
$newText = "Clustered only hotel locations only with k=130 and then assigned each Airbnb listing to a centroid if they were within 1km. 36% of Airbnb listings were disregarded for being too far away from a center."

$full = $tr.Text
$idx = $full.IndexOf($oldText)
if ($idx -lt 0) {
    throw "Target bullet text not found in shape '$($sh.Name)'"
}

$sub = $tr.Characters($idx + 1, $oldText.Length)
$sub.Text = $newText
